# enh: sanity check crit buckling
#
# Adds a "Final loads" sanity-check table (segments 10-18, Nxx/Nyy/Nzz plus
# normalised pnxx/pnyy/pnxy ratios and a concatenated summary string) below
# the existing optimisation log, plus an "average std" helper cell next to
# the pivot table on row 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# ---------------------------------------------------------------------
# 1. New "average std" helper next to the existing STDEV pivot (row 16)
# ---------------------------------------------------------------------
$ws.Range("Q16").Value = "average std"
$ws.Range("R16").Formula = "=AVERAGE(Q14:Y14)"
$ws.Range("Q14").Copy()
$ws.Range("R16").PasteSpecial(-4122)   # xlPasteFormats - reuse the existing 0.0 style

# ---------------------------------------------------------------------
# 2. "Final loads" sanity-check table
# ---------------------------------------------------------------------
$ws.Range("D65").Value = "Final loads"

$ws.Range("C66").Value = "panel"
$ws.Range("D66").Value = "Nxx"
$ws.Range("E66").Value = "Nyy"
$ws.Range("F66").Value = "Nzz"
$ws.Range("G66").Value = "pnxx"
$ws.Range("H66").Value = "pnyy"
$ws.Range("I66").Value = "pnxy"

# Raw loads per panel segment (10-18)
$ws.Range("C67").Value = 10
$ws.Range("C68").Value = 11
$ws.Range("C69").Value = 12
$ws.Range("C70").Value = 13
$ws.Range("C71").Value = 14
$ws.Range("C72").Value = 15
$ws.Range("C73").Value = 16
$ws.Range("C74").Value = 17
$ws.Range("C75").Value = 18

$ws.Range("D67").Value = -2413.2969069999999
$ws.Range("D68").Value = -1519.608782
$ws.Range("D69").Value = -449.3331657
$ws.Range("D70").Value = -448.05270000000002
$ws.Range("D71").Value = -1474.017625
$ws.Range("D72").Value = -2782.6806879999999
$ws.Range("D73").Value = -2831.1223749999999
$ws.Range("D74").Value = -1670.3654690000001
$ws.Range("D75").Value = -531.50832500000001

$ws.Range("E67").Value = -607.77428129999998
$ws.Range("E68").Value = -40.736502819999998
$ws.Range("E69").Value = 2.8955858750000001
$ws.Range("E70").Value = 30.941809379999999
$ws.Range("E71").Value = -155.11935
$ws.Range("E72").Value = -715.9270563
$ws.Range("E73").Value = -585.61413440000001
$ws.Range("E74").Value = -65.819476249999994
$ws.Range("E75").Value = 7.9851393750000002

$ws.Range("F67").Value = 291.44217190000001
$ws.Range("F68").Value = -29.929870000000001
$ws.Range("F69").Value = -61.942467499999999
$ws.Range("F70").Value = -176.95502819999999
$ws.Range("F71").Value = -209.19045
$ws.Range("F72").Value = -126.92502500000001
$ws.Range("F73").Value = 22.45525632
$ws.Range("F74").Value = -324.79375320000003
$ws.Range("F75").Value = -284.7068688

# Style: reuse the workbook's existing "segment header" style (C67:C75)
# and the existing "0" integer-display style (D:F 67:75) by pasting the
# formats from cells that already carry them.
$ws.Range("B62").Copy()
$ws.Range("C67:C75").PasteSpecial(-4122)

$ws.Range("Q5").Copy()
$ws.Range("D67:F75").PasteSpecial(-4122)

# Ratio formulas (|Nxx/Nxx|, |Nyy/Nxx|, |Nzz/Nxx|) and the concatenated
# "(g,h,i)" sanity-check string for every segment row.
for ($r = 67; $r -le 75; $r++) {
    $ws.Range("G$r").Formula = "=ABS(D$r/D$r)"
    $ws.Range("H$r").Formula = "=ABS(E$r/D$r)"
    $ws.Range("I$r").Formula = "=ABS(F$r/D$r)"
    $ws.Range("J$r").Formula = "=_xlfn.CONCAT(""("",G$r,"","",H$r,"","",I$r)"
}

# H/I columns get a new "0.00" display style - create it once, then reuse
# the resulting style for the rest of the block (and for H62/I62 below).
$ws.Range("H67").NumberFormat = "0.00"
$ws.Range("H67").Copy()
$ws.Range("H67:I75").PasteSpecial(-4122)
$ws.Range("H62").PasteSpecial(-4122)
$ws.Range("I62").PasteSpecial(-4122)

# Re-apply the formulas/values that the blanket format paste just
# overwrote on H67 (PasteSpecial(formats) does not touch H68:I75 values,
# but guard H67 explicitly so the first cell keeps its formula too).
$ws.Range("H67").Formula = "=ABS(E67/D67)"
$ws.Range("I67").Formula = "=ABS(F67/D67)"

# ---------------------------------------------------------------------
# 3. Manually transcribed (g,h,i) summary strings, pasted as plain text
# ---------------------------------------------------------------------
$ws.Range("I78").Value = "(0.25,0.12)"
$ws.Range("I79").Value = "(0.02,0.01)"
$ws.Range("I80").Value = "(0.01,0.14)"
$ws.Range("I81").Value = "(0.07,0.39)"
$ws.Range("I82").Value = "(0.11,0.14)"
$ws.Range("I83").Value = "(0.26,0.05)"
$ws.Range("I84").Value = "(0.20,0.01)"
$ws.Range("I85").Value = "(0.04,0.19)"
$ws.Range("I86").Value = "(0.01,0.53)"

# ---------------------------------------------------------------------
# 4. Cosmetic column-width tweaks that came along with the edit
# ---------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 8.1666666666667    # -> stored width 9
$ws.Columns.Item(18).ColumnWidth = 3.92               # -> stored width ~4.75 (splits off from 19:25)

Write-Output "Final loads sanity-check table written."
